# Scheduled price-refresh update.
#
# The upstream market-board price puller re-ran and produced new average
# sale prices for a subset of leve items. This script pushes the refreshed
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ values
# into each class/job sheet's table, which in turn changes the dependent
# LevePriceNQ/HQ and LeveProfitNQ/HQ figures for the affected rows.
#
# Sheets touched: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3490.5
$ws.Range("I6").Value = 3490.5
$ws.Range("K6").Value = 10471.5
$ws.Range("M6").Value = -10359.5
$ws.Range("H28").Value = 741456.0600000001
$ws.Range("I28").Value = 1111334.1
$ws.Range("J28").Value = 1700
$ws.Range("K28").Value = 1111334.1
$ws.Range("L28").Value = 1700
$ws.Range("M28").Value = -1110849.1
$ws.Range("N28").Value = -2670
$ws.Range("H112").Value = 9741240
$ws.Range("I112").Value = 1248.3334
$ws.Range("J112").Value = 12397602
$ws.Range("K112").Value = 3745.0002
$ws.Range("L112").Value = 37192806
$ws.Range("M112").Value = -2637.0002
$ws.Range("N112").Value = -37195022
$ws.Range("H129").Value = 1415.6818
$ws.Range("I129").Value = 597
$ws.Range("J129").Value = 1497.55
$ws.Range("K129").Value = 1791
$ws.Range("L129").Value = 4492.65
$ws.Range("M129").Value = 3209
$ws.Range("N129").Value = -14492.65
$ws.Range("H133").Value = 31434
$ws.Range("J133").Value = 31434
$ws.Range("L133").Value = 31434
$ws.Range("N133").Value = -41554
$ws.Range("H137").Value = 45455744
$ws.Range("I137").Value = 47620160
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 142860480
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -142857930
$ws.Range("N137").Value = -14100
$ws.Range("H138").Value = 4873381.5
$ws.Range("I138").Value = 2980646.8
$ws.Range("J138").Value = 5293989
$ws.Range("K138").Value = 8941940.399999999
$ws.Range("L138").Value = 15881967
$ws.Range("M138").Value = -8936800.399999999
$ws.Range("N138").Value = -15892247
$ws.Range("H141").Value = 1679.8064
$ws.Range("I141").Value = 1620.3214
$ws.Range("J141").Value = 2235
$ws.Range("K141").Value = 4860.9642
$ws.Range("L141").Value = 6705
$ws.Range("M141").Value = 319.0357999999997
$ws.Range("N141").Value = -17065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 900000000
$ws.Range("I10").Value = 900000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 900000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -899999830
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 17595.926
$ws.Range("I32").Value = 3037.1606
$ws.Range("J32").Value = 85536.836
$ws.Range("K32").Value = 3037.1606
$ws.Range("L32").Value = 85536.836
$ws.Range("M32").Value = -2750.1606
$ws.Range("N32").Value = -86110.836
$ws.Range("H45").Value = 900.6
$ws.Range("I45").Value = 867.5
$ws.Range("J45").Value = 950.25
$ws.Range("K45").Value = 867.5
$ws.Range("L45").Value = 950.25
$ws.Range("M45").Value = -490.5
$ws.Range("N45").Value = -1704.25
$ws.Range("H58").Value = 10600
$ws.Range("J58").Value = 10600
$ws.Range("L58").Value = 10600
$ws.Range("N58").Value = -11460
$ws.Range("H61").Value = 1220.3654
$ws.Range("I61").Value = 818.0732
$ws.Range("J61").Value = 2719.818
$ws.Range("K61").Value = 818.0732
$ws.Range("L61").Value = 2719.818
$ws.Range("M61").Value = -606.0732
$ws.Range("N61").Value = -3143.818
$ws.Range("H74").Value = 3105.8982
$ws.Range("I74").Value = 985.6667
$ws.Range("J74").Value = 8344.117
$ws.Range("K74").Value = 985.6667
$ws.Range("L74").Value = 8344.117
$ws.Range("M74").Value = -111.6667
$ws.Range("N74").Value = -10092.117
$ws.Range("H77").Value = 3105.8982
$ws.Range("I77").Value = 985.6667
$ws.Range("J77").Value = 8344.117
$ws.Range("K77").Value = 4928.3335
$ws.Range("L77").Value = 41720.585
$ws.Range("M77").Value = -560.3334999999997
$ws.Range("N77").Value = -50456.585
$ws.Range("H136").Value = 1220.3654
$ws.Range("I136").Value = 818.0732
$ws.Range("J136").Value = 2719.818
$ws.Range("K136").Value = 2454.2196
$ws.Range("L136").Value = 8159.454000000001
$ws.Range("M136").Value = 95.78039999999964
$ws.Range("N136").Value = -13259.454
$ws.Range("H139").Value = 44500
$ws.Range("J139").Value = 44500
$ws.Range("L139").Value = 44500
$ws.Range("N139").Value = -54780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 41440
$ws.Range("J59").Value = 46800
$ws.Range("L59").Value = 46800
$ws.Range("N59").Value = -48494
$ws.Range("H94").Value = 1321.9333
$ws.Range("I94").Value = 1462.5
$ws.Range("J94").Value = 1161.2858
$ws.Range("K94").Value = 1462.5
$ws.Range("L94").Value = 1161.2858
$ws.Range("M94").Value = -1011.5
$ws.Range("N94").Value = -2063.2858
$ws.Range("H132").Value = 26293.334
$ws.Range("J132").Value = 26293.334
$ws.Range("L132").Value = 26293.334
$ws.Range("N132").Value = -36413.334
$ws.Range("H133").Value = 49890
$ws.Range("J133").Value = 49890
$ws.Range("L133").Value = 49890
$ws.Range("N133").Value = -60010
$ws.Range("H134").Value = 19610480
$ws.Range("I134").Value = 31251362
$ws.Range("J134").Value = 4784.9473
$ws.Range("K134").Value = 93754086
$ws.Range("L134").Value = 14354.8419
$ws.Range("M134").Value = -93751551
$ws.Range("N134").Value = -19424.8419

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1483.4186
$ws.Range("I31").Value = 948.6774
$ws.Range("K31").Value = 948.6774
$ws.Range("M31").Value = -653.6774
$ws.Range("H34").Value = 1483.4186
$ws.Range("I34").Value = 948.6774
$ws.Range("K34").Value = 948.6774
$ws.Range("M34").Value = -746.6774
$ws.Range("H58").Value = 2743.5334
$ws.Range("I58").Value = 1494.4
$ws.Range("J58").Value = 3368.1
$ws.Range("K58").Value = 1494.4
$ws.Range("L58").Value = 3368.1
$ws.Range("M58").Value = -1291.4
$ws.Range("N58").Value = -3774.1
$ws.Range("H132").Value = 3417.2942
$ws.Range("I132").Value = 2495.111
$ws.Range("J132").Value = 4454.75
$ws.Range("K132").Value = 7485.333
$ws.Range("L132").Value = 13364.25
$ws.Range("M132").Value = -4955.333
$ws.Range("N132").Value = -18424.25
$ws.Range("H136").Value = 2743.5334
$ws.Range("I136").Value = 1494.4
$ws.Range("J136").Value = 3368.1
$ws.Range("K136").Value = 4483.200000000001
$ws.Range("L136").Value = 10104.3
$ws.Range("M136").Value = -1933.200000000001
$ws.Range("N136").Value = -15204.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1308.6666
$ws.Range("I5").Value = 517
$ws.Range("J5").Value = 2043.7858
$ws.Range("K5").Value = 1551
$ws.Range("L5").Value = 6131.357400000001
$ws.Range("M5").Value = -1439
$ws.Range("N5").Value = -6355.357400000001
$ws.Range("H60").Value = 1695.8695
$ws.Range("I60").Value = 130
$ws.Range("J60").Value = 2248.5293
$ws.Range("K60").Value = 390
$ws.Range("L60").Value = 6745.5879
$ws.Range("M60").Value = -139
$ws.Range("N60").Value = -7247.5879
$ws.Range("H122").Value = 845.5833
$ws.Range("I122").Value = 288
$ws.Range("K122").Value = 2592
$ws.Range("M122").Value = -142
$ws.Range("H131").Value = 1530.7142
$ws.Range("I131").Value = 501.42856
$ws.Range("J131").Value = 1702.262
$ws.Range("K131").Value = 1504.28568
$ws.Range("L131").Value = 5106.786
$ws.Range("M131").Value = 3535.71432
$ws.Range("N131").Value = -15186.786
$ws.Range("H135").Value = 1308.6666
$ws.Range("I135").Value = 517
$ws.Range("J135").Value = 2043.7858
$ws.Range("K135").Value = 4653
$ws.Range("L135").Value = 18394.0722
$ws.Range("M135").Value = -2118
$ws.Range("N135").Value = -23464.0722
$ws.Range("H136").Value = 2495.6155
$ws.Range("I136").Value = 1657.1428
$ws.Range("J136").Value = 2804.5264
$ws.Range("K136").Value = 4971.428400000001
$ws.Range("L136").Value = 8413.5792
$ws.Range("M136").Value = 128.5715999999993
$ws.Range("N136").Value = -18613.5792
$ws.Range("H137").Value = 4594795
$ws.Range("I137").Value = 8336121
$ws.Range("J137").Value = 105203.3
$ws.Range("K137").Value = 25008363
$ws.Range("L137").Value = 315609.9
$ws.Range("M137").Value = -25003263
$ws.Range("N137").Value = -325809.9
$ws.Range("H139").Value = 2005.4615
$ws.Range("I139").Value = 1787.6765
$ws.Range("K139").Value = 5363.029500000001
$ws.Range("M139").Value = -223.0295000000006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1011226.8
$ws.Range("I122").Value = 1112261.5
$ws.Range("K122").Value = 3336784.5
$ws.Range("M122").Value = -3334334.5
$ws.Range("H132").Value = 3406.628
$ws.Range("I132").Value = 3043.5
$ws.Range("J132").Value = 4463
$ws.Range("K132").Value = 9130.5
$ws.Range("L132").Value = 13389
$ws.Range("M132").Value = -6600.5
$ws.Range("N132").Value = -18449
$ws.Range("H138").Value = 64000
$ws.Range("J138").Value = 64000
$ws.Range("L138").Value = 64000
$ws.Range("N138").Value = -74280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3205.5557
$ws.Range("I7").Value = 2050
$ws.Range("J7").Value = 3535.7144
$ws.Range("K7").Value = 2050
$ws.Range("L7").Value = 3535.7144
$ws.Range("M7").Value = -1938
$ws.Range("N7").Value = -3759.7144
$ws.Range("H126").Value = 3205.5557
$ws.Range("I126").Value = 2050
$ws.Range("J126").Value = 3535.7144
$ws.Range("K126").Value = 6150
$ws.Range("L126").Value = 10607.1432
$ws.Range("M126").Value = -3680
$ws.Range("N126").Value = -15547.1432
$ws.Range("H132").Value = 4607.1333
$ws.Range("I132").Value = 4062.8845
$ws.Range("J132").Value = 5351.8945
$ws.Range("K132").Value = 12188.6535
$ws.Range("L132").Value = 16055.6835
$ws.Range("M132").Value = -9658.6535
$ws.Range("N132").Value = -21115.6835
$ws.Range("H136").Value = 3202.7637
$ws.Range("I136").Value = 1780.6522
$ws.Range("J136").Value = 10471.333
$ws.Range("K136").Value = 5341.9566
$ws.Range("L136").Value = 31413.999
$ws.Range("M136").Value = -2791.9566
$ws.Range("N136").Value = -36513.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3338
$ws.Range("H107").Value = 4274893
$ws.Range("I107").Value = 6945694.5
$ws.Range("J107").Value = 1610
$ws.Range("K107").Value = 20837083.5
$ws.Range("L107").Value = 4830
$ws.Range("M107").Value = -20835163.5
$ws.Range("N107").Value = -8670
$ws.Range("H126").Value = 143675.42
$ws.Range("I126").Value = 200775.6
$ws.Range("J126").Value = 925
$ws.Range("K126").Value = 602326.8
$ws.Range("L126").Value = 2775
$ws.Range("M126").Value = -599856.8
$ws.Range("N126").Value = -7715
$ws.Range("H132").Value = 20003120
$ws.Range("I132").Value = 29414432
$ws.Range("K132").Value = 88243296
$ws.Range("M132").Value = -88240766
$ws.Range("H136").Value = 6667912
$ws.Range("I136").Value = 7937299.5
$ws.Range("J136").Value = 3630.125
$ws.Range("K136").Value = 23811898.5
$ws.Range("L136").Value = 10890.375
$ws.Range("M136").Value = -23809348.5
$ws.Range("N136").Value = -15990.375
